# Update the Bacillus tree distance-comparison matrix on Sheet1 so the script
# reflects a full, completed run (NC_017188 vs NC_017190 tree comparison).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell address -> new (text) value. Values are written as text so they
# keep their "N.N" formatting (matching the rest of the matrix) instead of
# being re-interpreted as numbers.
$cellUpdates = [ordered]@{
    "F4" = "5.0"
    "I4" = "5.0"
    "J4" = "7.0"
    "K4" = "5.0"
    "R4" = "5.0"
    "S4" = "22.0"
    "V4" = "5.0"
    "C5" = "5.0"
    "E5" = "5.0"
    "G5" = "3.0"
    "H5" = "3.0"
    "J5" = "6.0"
    "K5" = "6.0"
    "L5" = "18.0"
    "R5" = "2.0"
    "V5" = "3.0"
    "F6" = "3.0"
    "I6" = "5.0"
    "J6" = "7.0"
    "R6" = "3.0"
    "V6" = "4.0"
    "F7" = "3.0"
    "I7" = "5.0"
    "J7" = "7.0"
    "R7" = "3.0"
    "V7" = "4.0"
    "C8" = "5.0"
    "E8" = "5.0"
    "G8" = "5.0"
    "H8" = "5.0"
    "J8" = "6.0"
    "K8" = "7.0"
    "R8" = "5.0"
    "S8" = "22.0"
    "V8" = "5.0"
    "C9" = "7.0"
    "E9" = "7.0"
    "F9" = "6.0"
    "G9" = "7.0"
    "H9" = "7.0"
    "I9" = "6.0"
    "K9" = "7.0"
    "L9" = "14.5"
    "R9" = "7.0"
    "S9" = "22.0"
    "V9" = "6.0"
    "C10" = "5.0"
    "E10" = "5.0"
    "F10" = "6.0"
    "I10" = "7.0"
    "J10" = "7.0"
    "R10" = "6.5"
    "V10" = "7.0"
    "F11" = "18.0"
    "J11" = "14.5"
    "R11" = "19.0"
    "S11" = "19.0"
    "V11" = "16.0"
    "C16" = "5.0"
    "E16" = "5.0"
    "F16" = "2.0"
    "G16" = "3.0"
    "H16" = "3.0"
    "I16" = "5.0"
    "J16" = "7.0"
    "K16" = "6.5"
    "L16" = "19.0"
    "S16" = "23.5"
    "V16" = "4.0"
    "C17" = "22.0"
    "E17" = "22.0"
    "I17" = "22.0"
    "J17" = "22.0"
    "L17" = "19.0"
    "R17" = "23.5"
    "V17" = "21.0"
    "C20" = "5.0"
    "E20" = "5.0"
    "F20" = "3.0"
    "G20" = "4.0"
    "H20" = "4.0"
    "I20" = "5.0"
    "J20" = "6.0"
    "K20" = "7.0"
    "L20" = "16.0"
    "R20" = "4.0"
    "S20" = "21.0"
}

foreach ($addr in $cellUpdates.Keys) {
    $cell = $ws.Range($addr)
    # Force text format first so Excel stores the value as a string (matching
    # the workbook's existing convention of storing these numbers as text),
    # then restore the default "Normal" style so no stray formatting is left
    # behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $cellUpdates[$addr]
    $cell.Style = "Normal"
}
